$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$ws1.Range("G2").Value = "2016-10-13 13:42:18"
$ws3.Range("H2").Value = "2016-10-13 13:42:18"

$ws2.Range("H2").Value = "2016-10-13 13:42:07"
$ws2.Range("K2").Value = "2016-10-13 13:42:45"

$ws3.Range("K2").Value = "2016-10-13 13:43:01"
